$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.174.63'
$ws.Range("E2").Value = '  +1.01%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.081.53'
$ws.Range("E3").Value = '  +0.30%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("E4").Value = '  +0.02%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.17'
$ws.Range("E5").Value = '  +1.13%  '

# Row 6: 'Solana' -> 'Solana'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.13'
$ws.Range("E6").Value = '  +0.91%  '

# Row 7: 'USDC' -> 'USDC'
$ws.Range("E7").Value = '  +0.18%  '

# Row 8: 'LidoStakedEther' -> 'LidoStakedEther'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.078.47'
$ws.Range("E8").Value = '  +0.36%  '

# Row 9: 'XRP' -> 'XRP'
$ws.Range("E9").Value = '  +0.85%  '

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range("E10").Value = '  +1.82%  '

# Row 11: 'Toncoin' -> 'Toncoin'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.12'
$ws.Range("E11").Value = '  -6.06%  '

# Row 12: 'Cardano' -> 'Cardano'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.473'
$ws.Range("E12").Value = '  +3.05%  '

# Row 13: 'ShibaInu' -> 'ShibaInu'
$ws.Range("E13").Value = '  +0.27%  '

# Row 14: 'Avalanche' -> 'Avalanche'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.06'
$ws.Range("E14").Value = '  +0.16%  '

# Row 15: 'WrappedliquidstakedEther2.0' -> 'WrappedliquidstakedEther2.0'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.597.04'
$ws.Range("E15").Value = '  +0.91%  '

# Row 16: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.196.57'
$ws.Range("E16").Value = '  +1.09%  '

# Row 17: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.086.17'
$ws.Range("E17").Value = '  +0.48%  '

# Row 18: 'TRON' -> 'TRON'
$ws.Range("E18").Value = '  +1.15%  '

# Row 19: 'Polkadot' -> 'Polkadot'
$ws.Range("E19").Value = '  -1.06%  '

# Row 20: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.70'
$ws.Range("E20").Value = '  -1.03%  '

# Row 21: 'Chainlink' -> 'Chainlink'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.01'
$ws.Range("E21").Value = '  +0.62%  '

# Row 22: 'Polygon' -> 'Polygon'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.678'

# Row 23: 'Uniswap' -> 'Uniswap'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("E23").Value = '  +3.24%  '

# Row 24: 'InternetComputer(DFINITY)' -> 'InternetComputer(DFINITY)'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.08'
$ws.Range("E24").Value = '  +10.10%  '

# Row 25: 'Litecoin' -> 'Litecoin'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.39'
$ws.Range("E25").Value = '  +0.23%  '

# Row 26: 'Dai' -> 'Dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.03%  '

# Row 27: 'PancakeSwap' -> 'PancakeSwap'
$ws.Range("E27").Value = '  +0.51%  '

# Row 28: 'RenderToken' -> 'RenderToken'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.98'
$ws.Range("E28").Value = '  +0.63%  '

# Row 29: 'ImmutableX' -> 'ImmutableX'
$ws.Range("E29").Value = '  +1.43%  '

# Row 30: 'FirstDigitalUSD' -> 'FirstDigitalUSD'
$ws.Range("E30").Value = '  +0.04%  '

# Row 31: 'EthereumClassic' -> 'EthereumClassic'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.29'

# Row 32: 'Mantle' -> 'Mantle'
$ws.Range("E32").Value = '  -1.65%  '

# Row 33: 'Stacks' -> 'Stacks'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.47'
$ws.Range("E33").Value = '  +0.45%  '

# Row 34: 'NEARProtocol' -> 'NEARProtocol'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.64'
$ws.Range("E34").Value = '  -1.16%  '

# Row 35: 'Filecoin' -> 'Filecoin'
$ws.Range("E35").Value = '  +2.95%  '

# Row 36: 'OKB' -> 'OKB'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.35'
$ws.Range("E36").Value = '  -0.46%  '

# Row 37: 'dogwifhat' -> 'VeChain'
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0407'
$ws.Range("E37").Value = '  +1.81%  '

# Row 38: 'VeChain' -> 'dogwifhat'
$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.95'
$ws.Range("E38").Value = '  +13.79%  '

# Row 39: 'Bittensor' -> 'Bittensor'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '438.42'
$ws.Range("E39").Value = '  -6.39%  '

# Row 40: 'Hedera' -> 'Hedera'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0810'
$ws.Range("E40").Value = '  -1.98%  '

# Row 41: 'Maker' -> 'Maker'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.957.62'
$ws.Range("E41").Value = '  -2.95%  '

# Row 42: 'Cosmos' -> 'Cosmos'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.20'
$ws.Range("E42").Value = '  -0.66%  '

# Row 43: 'Kaspa' -> 'Kaspa'
$ws.Range("E43").Value = '  -4.84%  '

# Row 44: 'InjectiveProtocol' -> 'InjectiveProtocol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.17'
$ws.Range("E44").Value = '  +1.31%  '

# Row 45: 'TheGraph' -> 'TheGraph'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.260'
$ws.Range("E45").Value = '  +1.30%  '

# Row 47: 'Fetch.AI' -> 'Fetch.AI'
$ws.Range("E47").Value = '  +3.34%  '

# Row 48: 'Stellar' -> 'Stellar'
$ws.Range("E48").Value = '  +0.90%  '

# Row 49: 'Monero' -> 'Monero'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '118.11'
$ws.Range("E49").Value = '  +1.12%  '

# Row 50: 'PEPE' -> 'PEPE'
$ws.Range("E50").Value = '  +0.62%  '

# Row 51: 'ThetaToken' -> 'ThetaToken'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.08'
$ws.Range("E51").Value = '  -0.60%  '
